# Insert a new price record row at row 85 (Flame Seedless, Provincia de Copiapó)
# Existing rows 85-165 shift down to 86-166.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85; existing row 85 and below shift down one row.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the same constant columns used
# throughout the sheet (A,B,C,E,F,G,H,I,J,L), plus the new record's values.
$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 45271
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100109
$ws.Cells.Item(85, 8).Value = "Uva"
$ws.Cells.Item(85, 9).Value = 100109001
$ws.Cells.Item(85, 10).Value = "Uva"
$ws.Cells.Item(85, 11).Value = "Flame Seedless"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 700
$ws.Cells.Item(85, 14).Value = 12000
$ws.Cells.Item(85, 15).Value = 13000
$ws.Cells.Item(85, 16).Value = 12500
$ws.Cells.Item(85, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Copiapó"
$ws.Cells.Item(85, 19).Value = 1250
$ws.Cells.Item(85, 20).Value = 10
